# Apply the balance-table extension (rows 62-102) to Hoja1 (sheet1),
# style the A1 header cell, and finish with the J96 selection — mirroring
# a user who dragged the fill handle on A61:H61 down to row 102 and then
# painted the "Bueno" style onto A1 to match the rest of the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastOld = 61          # last pre-existing data row
$lastNew = 102          # last row after the extension
$firstBlockEnd = 84      # end of the first fill-down block (E/F/G/H shared formulas)

# --- Column A: round index, simply increases by one each row ------------
for ($r = $lastOld + 1; $r -le $lastNew; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# --- Columns B & C: each row keeps its own (non-shared) formula that ----
# references column A of the row directly above it, matching rows 3-61.
for ($r = $lastOld + 1; $r -le $lastNew; $r++) {
    $prevRow = $r - 1
    $ws.Cells.Item($r, 2).Formula = "= `$B`$3 * hits_1_round_behind ^ A$prevRow"
    $ws.Cells.Item($r, 3).Formula = "= `$B`$3*hits_1_round_behind^A$prevRow*brute_life_multiplier"
}

# --- Columns E, F, G, H: filled as two blocks (62-84, 85-102) so each ---
# block becomes one Excel shared formula, exactly as in the source file.
$blocks = @(
    @{ Start = $lastOld + 1; End = $firstBlockEnd },
    @{ Start = $firstBlockEnd + 1; End = $lastNew }
)

foreach ($block in $blocks) {
    $s = $block.Start
    $e = $block.End
    $prevRow = $s - 1

    $ws.Range("E$s`:E$e").Formula = "= base_damage * hits_1_round_behind ^ A$prevRow"
    $ws.Range("F$s`:F$e").Formula = "=ROUNDUP(B$s/base_damage, 0)"
    $ws.Range("G$s`:G$e").Formula = "=E$s/2"
    $ws.Range("H$s`:H$e").Formula = "=ROUNDUP(B$s/explosion_shot_base, 0)"
}

# --- Header cell A1 now gets the same "Bueno" style as its neighbours ---
$ws.Range("A1").Style = "Bueno"

# --- Final selection left at J96, where the user ended up -------------
$ws.Range("J96").Select() | Out-Null
